$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update period labels (E16:E20) to reflect the new reversed chronological order
$ws.Range("E16").Value = "2002"
$ws.Range("E17").Value = "2001"
$ws.Range("E18").Value = "1912"
$ws.Range("E19").Value = "1911"
$ws.Range("E20").Value = "1910"

# Update Valor Mora values (F16:F20) - values swap so 1910 row now has 8833, others 33125
$ws.Range("F16").Value = 33125
$ws.Range("F17").Value = 33125
$ws.Range("F18").Value = 33125
$ws.Range("F19").Value = 33125
$ws.Range("F20").Value = 8833
